{"js": "// Update the worksheet date stamp and the twenty-five \"three-digit \u00d7\n// one-digit\" multiplication problems/answers embedded in the table cells.\n// Each old value is unique within the document, so a scoped search/replace\n// per pair is unambiguous and order-independent.\nconst replacements = [\n  [\"2025-03-28 Friday\", \"2025-03-29 Saturday\"],\n  [\"479\u00d72=958\", \"312\u00d72=624\"],\n  [\"715\u00d74=2860\", \"356\u00d74=1424\"],\n  [\"481\u00d73=1443\", \"394\u00d76=2364\"],\n  [\"671\u00d72=1342\", \"409\u00d74=1636\"],\n  [\"972\u00d77=6804\", \"426\u00d72=852\"],\n  [\"285\u00d77=1995\", \"490\u00d78=3920\"],\n  [\"403\u00d76=2418\", \"500\u00d79=4500\"],\n  [\"736\u00d75=3680\", \"684\u00d73=2052\"],\n  [\"412\u00d72=824\", \"476\u00d74=1904\"],\n  [\"126\u00d74=504\", \"300\u00d73=900\"],\n  [\"139\u00d73=417\", \"274\u00d77=1918\"],\n  [\"528\u00d76=3168\", \"187\u00d75=935\"],\n  [\"806\u00d79=7254\", \"930\u00d74=3720\"],\n  [\"102\u00d78=816\", \"905\u00d79=8145\"],\n  [\"721\u00d74=2884\", \"954\u00d76=5724\"],\n  [\"385\u00d74=1540\", \"797\u00d72=1594\"],\n  [\"519\u00d79=4671\", \"513\u00d75=2565\"],\n  [\"927\u00d78=7416\", \"400\u00d73=1200\"],\n  [\"320\u00d73=960\", \"586\u00d78=4688\"],\n  [\"249\u00d72=498\", \"941\u00d73=2823\"],\n  [\"631\u00d79=5679\", \"536\u00d74=2144\"],\n  [\"517\u00d74=2068\", \"392\u00d77=2744\"],\n  [\"550\u00d79=4950\", \"741\u00d73=2223\"],\n  [\"353\u00d74=1412\", \"842\u00d77=5894\"],\n  [\"770\u00d74=3080\", \"346\u00d77=2422\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for ${JSON.stringify(oldText)}, found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date stamp and the twenty-five \"three-digit \u00d7\n# one-digit\" multiplication problems/answers embedded in the table cells.\n# Each old value is unique within the document, so a simple Find/Replace\n# (wdReplaceAll, scoped to the whole document) per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-28 Friday\", \"2025-03-29 Saturday\"),\n    @(\"479\u00d72=958\", \"312\u00d72=624\"),\n    @(\"715\u00d74=2860\", \"356\u00d74=1424\"),\n    @(\"481\u00d73=1443\", \"394\u00d76=2364\"),\n    @(\"671\u00d72=1342\", \"409\u00d74=1636\"),\n    @(\"972\u00d77=6804\", \"426\u00d72=852\"),\n    @(\"285\u00d77=1995\", \"490\u00d78=3920\"),\n    @(\"403\u00d76=2418\", \"500\u00d79=4500\"),\n    @(\"736\u00d75=3680\", \"684\u00d73=2052\"),\n    @(\"412\u00d72=824\", \"476\u00d74=1904\"),\n    @(\"126\u00d74=504\", \"300\u00d73=900\"),\n    @(\"139\u00d73=417\", \"274\u00d77=1918\"),\n    @(\"528\u00d76=3168\", \"187\u00d75=935\"),\n    @(\"806\u00d79=7254\", \"930\u00d74=3720\"),\n    @(\"102\u00d78=816\", \"905\u00d79=8145\"),\n    @(\"721\u00d74=2884\", \"954\u00d76=5724\"),\n    @(\"385\u00d74=1540\", \"797\u00d72=1594\"),\n    @(\"519\u00d79=4671\", \"513\u00d75=2565\"),\n    @(\"927\u00d78=7416\", \"400\u00d73=1200\"),\n    @(\"320\u00d73=960\", \"586\u00d78=4688\"),\n    @(\"249\u00d72=498\", \"941\u00d73=2823\"),\n    @(\"631\u00d79=5679\", \"536\u00d74=2144\"),\n    @(\"517\u00d74=2068\", \"392\u00d77=2744\"),\n    @(\"550\u00d79=4950\", \"741\u00d73=2223\"),\n    @(\"353\u00d74=1412\", \"842\u00d77=5894\"),\n    @(\"770\u00d74=3080\", \"346\u00d77=2422\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $oldText\n    $range.Find.Replacement.Text = $newText\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Find/Replace failed: could not locate '$oldText' in the document.\"\n    }\n}\n"}
